$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.359.58"

$ws.Range("D3").Value = "1.867.56"
$ws.Range("E3").Value = "  +0.30%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.82%  "

$ws.Range("E6").Value = "  +0.02%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4714"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.64%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2866"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.11%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.41"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.02%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07877"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.72%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "96.98"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.17%  "

$ws.Range("D13").Value = "1.865.69"
$ws.Range("E13").Value = "  +0.14%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6914"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.59%  "

$ws.Range("E15").Value = "  -1.15%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "268.59"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.40%  "

$ws.Range("D17").Value = "30.326.42"
$ws.Range("E17").Value = "  +0.37%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.93"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.13%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007661"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.74%  "

$ws.Range("E20").Value = "  -0.07%  "

$ws.Range("D21").Value = "2.115.79"
$ws.Range("E21").Value = "  +0.09%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.03%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.230"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.66%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.183"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.02%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.402"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.07%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.37"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.47%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.89"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.06%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.945"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.66%  "

$ws.Range("E29").Value = "  -2.04%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09920"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.82%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.376"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.15%  "

$ws.Range("E32").Value = "  -0.70%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.057"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.14%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04749"
$ws.Range("D34").Style = "Normal"

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.135"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.42%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7030"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.03%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.720"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.49%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01876"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.32%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.793"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.90%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.307"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.21%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "73.33"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.26%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.953"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.86%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8437"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.25%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4178"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.40%  "

$ws.Range("E45").Value = "  +0.10%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "103.06"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.13%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "964.49"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.85%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.115"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.52%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.143"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.14%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.52"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.19%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05680"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.50%  "
